$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 15: TP11 entry (mirrors the TP10 row 11/12 comment pattern)
$ws.Range("A15").Value = "TP11"
$ws.Range("B15").Value = "Visualização arquitetural (dois apresentaram)
- antlr4, json, java, mayavi, networkx
- análise coleta tudo em json e plotar
- grafo, mas com cor para cada componente.
- Fizeram algo buscando um code city.
- Atividade desempenhada de forma como esperada."
$ws.Range("C15").Value = 7.5

# Copy formatting from row 11 (same visual pattern: A=s3, B=s17 wrap, C=s13)
$ws.Range("A11").Copy()
$ws.Range("A15").PasteSpecial(-4122) | Out-Null
$ws.Range("B11").Copy()
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("C11").Copy()
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row height for the new wrapped comment row
$ws.Rows.Item(15).RowHeight = 102

# Update the view: scroll position & active selection
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E13").Select()
